# feat: new final map
# Add a new row of reference data to the table and update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new reference row (row 3): Author, Title, Lincese, URL
$ws.Range("A3").Value = "ChatGPT 4o"
$ws.Range("B3").Value = "Pixel Dirt Path"
$ws.Range("C3").Value = "NULL"
$ws.Range("D3").Value = "NULL"

# Update the selected cell to match the author's final cursor position
$ws.Range("D5").Select()
